$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.695.58"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.01"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.36"
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4684"
$ws.Range("E7").Value = "  +4.05%  "
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.00"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08019"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.022"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.18"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.926"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.120"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001047"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.51"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06621"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.21"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.694.84"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.96"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.312"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.070.27"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.32"
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.082"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.546"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.53"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9679"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09484"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.446"
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.600"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.305"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02255"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06057"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.229"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.120"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1889"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.23"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.256"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5689"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.16"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.389"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.931"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06836"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.81"
$ws.Range("E51").Value = "  +5.00%  "
